# Auto-generated edit script: updates LevePrice/LeveProfit metrics for Behemoth_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2089.6
$ws.Range("J70").Value = 1450
$ws.Range("L70").Value = 4350
$ws.Range("N70").Value = -4890
$ws.Range("H73").Value = 2089.6
$ws.Range("J73").Value = 1450
$ws.Range("L73").Value = 4350
$ws.Range("N73").Value = -6222
$ws.Range("H80").Value = 924
$ws.Range("J80").Value = 1131.6666
$ws.Range("L80").Value = 3394.9998
$ws.Range("N80").Value = -5390.9998
$ws.Range("H83").Value = 924
$ws.Range("J83").Value = 1131.6666
$ws.Range("L83").Value = 10184.9994
$ws.Range("N83").Value = -20168.9994
$ws.Range("H100").Value = 1598
$ws.Range("I100").Value = 1531
$ws.Range("K100").Value = 1531
$ws.Range("M100").Value = -990

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 352.66666
$ws.Range("I5").Value = 219.8
$ws.Range("K5").Value = 219.8
$ws.Range("M5").Value = -107.8
$ws.Range("H32").Value = 13890841
$ws.Range("I32").Value = 13890841
$ws.Range("K32").Value = 13890841
$ws.Range("M32").Value = -13890554
$ws.Range("H63").Value = 3027.7778
$ws.Range("I63").Value = 3290.6667
$ws.Range("J63").Value = 2502
$ws.Range("K63").Value = 3290.6667
$ws.Range("L63").Value = 2502
$ws.Range("M63").Value = -2604.6667
$ws.Range("N63").Value = -3874
$ws.Range("H66").Value = 3027.7778
$ws.Range("I66").Value = 3290.6667
$ws.Range("J66").Value = 2502
$ws.Range("K66").Value = 16453.3335
$ws.Range("L66").Value = 12510
$ws.Range("M66").Value = -13021.3335
$ws.Range("N66").Value = -19374
$ws.Range("H74").Value = 7150328.5
$ws.Range("I74").Value = 11907020
$ws.Range("K74").Value = 11907020
$ws.Range("M74").Value = -11906146
$ws.Range("H77").Value = 7150328.5
$ws.Range("I77").Value = 11907020
$ws.Range("K77").Value = 59535100
$ws.Range("M77").Value = -59530732
$ws.Range("H95").Value = 33902.668
$ws.Range("J95").Value = 33902.668
$ws.Range("L95").Value = 33902.668
$ws.Range("N95").Value = -39394.668
$ws.Range("H113").Value = 109995
$ws.Range("J113").Value = 109995
$ws.Range("L113").Value = 109995
$ws.Range("N113").Value = -118673
$ws.Range("H120").Value = 111000
$ws.Range("J120").Value = 111000
$ws.Range("L120").Value = 111000
$ws.Range("N120").Value = -120676
$ws.Range("H139").Value = 59998.75
$ws.Range("J139").Value = 59998.75
$ws.Range("L139").Value = 59998.75
$ws.Range("N139").Value = -70278.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 352.66666
$ws.Range("I4").Value = 219.8
$ws.Range("K4").Value = 219.8
$ws.Range("M4").Value = -104.8
$ws.Range("H86").Value = 3180.111
$ws.Range("I86").Value = 3031.7144
$ws.Range("K86").Value = 3031.7144
$ws.Range("M86").Value = -1908.7144
$ws.Range("H89").Value = 3180.111
$ws.Range("I89").Value = 3031.7144
$ws.Range("K89").Value = 15158.572
$ws.Range("M89").Value = -9542.572
$ws.Range("H94").Value = 1039.2858
$ws.Range("I94").Value = 1041.5385
$ws.Range("K94").Value = 1041.5385
$ws.Range("M94").Value = -590.5385000000001
$ws.Range("H107").Value = 1293.8235
$ws.Range("I107").Value = 1285.0667
$ws.Range("K107").Value = 1285.0667
$ws.Range("M107").Value = 634.9332999999999
$ws.Range("H134").Value = 51594.547
$ws.Range("I134").Value = 2815.75
$ws.Range("K134").Value = 8447.25
$ws.Range("M134").Value = -5912.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 617988.7
$ws.Range("I31").Value = 4061.35
$ws.Range("K31").Value = 4061.35
$ws.Range("M31").Value = -3766.35
$ws.Range("H34").Value = 617988.7
$ws.Range("I34").Value = 4061.35
$ws.Range("K34").Value = 4061.35
$ws.Range("M34").Value = -3859.35
$ws.Range("H99").Value = 3791.1
$ws.Range("I99").Value = 3554.5
$ws.Range("J99").Value = 4737.5
$ws.Range("K99").Value = 3554.5
$ws.Range("L99").Value = 4737.5
$ws.Range("M99").Value = -2056.5
$ws.Range("N99").Value = -7733.5
$ws.Range("H117").Value = 82001
$ws.Range("J117").Value = 82001
$ws.Range("L117").Value = 82001
$ws.Range("N117").Value = -91179
$ws.Range("H125").Value = 251587
$ws.Range("J125").Value = 251587
$ws.Range("L125").Value = 251587
$ws.Range("N125").Value = -256507
$ws.Range("H126").Value = 3791.1
$ws.Range("I126").Value = 3554.5
$ws.Range("J126").Value = 4737.5
$ws.Range("K126").Value = 10663.5
$ws.Range("L126").Value = 14212.5
$ws.Range("M126").Value = -8193.5
$ws.Range("N126").Value = -19152.5
$ws.Range("H132").Value = 2847.2727
$ws.Range("I132").Value = 2832.1
$ws.Range("K132").Value = 8496.299999999999
$ws.Range("M132").Value = -5966.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1052.125
$ws.Range("J46").Value = 1997.25
$ws.Range("L46").Value = 5991.75
$ws.Range("N46").Value = -6173.75
$ws.Range("H51").Value = 16718.908
$ws.Range("I51").Value = 9856.143
$ws.Range("K51").Value = 29568.429
$ws.Range("M51").Value = -29108.429
$ws.Range("H81").Value = 4721.5
$ws.Range("J81").Value = 4444
$ws.Range("L81").Value = 13332
$ws.Range("N81").Value = -15578
$ws.Range("H84").Value = 4721.5
$ws.Range("J84").Value = 4444
$ws.Range("L84").Value = 39996
$ws.Range("N84").Value = -51228
$ws.Range("H129").Value = 23881202
$ws.Range("I129").Value = 4166
$ws.Range("J129").Value = 66859868
$ws.Range("K129").Value = 12498
$ws.Range("L129").Value = 200579604
$ws.Range("M129").Value = -7498
$ws.Range("N129").Value = -200589604

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9216751
$ws.Range("I11").Value = 9294000
$ws.Range("K11").Value = 9294000
$ws.Range("M11").Value = -9293861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 55556004
$ws.Range("I55").Value = 83333770
$ws.Range("J55").Value = 467
$ws.Range("K55").Value = 83333770
$ws.Range("L55").Value = 467
$ws.Range("M55").Value = -83333597
$ws.Range("N55").Value = -813
$ws.Range("H68").Value = 2580.125
$ws.Range("I68").Value = 2571.4546
$ws.Range("K68").Value = 2571.4546
$ws.Range("M68").Value = -1822.4546
$ws.Range("H71").Value = 2580.125
$ws.Range("I71").Value = 2571.4546
$ws.Range("K71").Value = 12857.273
$ws.Range("M71").Value = -9113.273000000001
$ws.Range("H82").Value = 1181
$ws.Range("J82").Value = 2000.5
$ws.Range("L82").Value = 2000.5
$ws.Range("N82").Value = -2722.5
$ws.Range("H85").Value = 1181
$ws.Range("J85").Value = 2000.5
$ws.Range("L85").Value = 2000.5
$ws.Range("N85").Value = -4496.5
$ws.Range("H122").Value = 7094.8335
$ws.Range("I122").Value = 5291.6665
$ws.Range("K122").Value = 15874.9995
$ws.Range("M122").Value = -13424.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = 6
$ws.Range("K9").Value = 6
$ws.Range("M9").Value = 134
$ws.Range("H17").Value = 337999.66
$ws.Range("I17").Value = 337999.66
$ws.Range("K17").Value = 337999.66
$ws.Range("M17").Value = -337827.66
$ws.Range("H62").Value = 4088321
$ws.Range("I62").Value = 7391.225
$ws.Range("K62").Value = 7391.225
$ws.Range("M62").Value = -6767.225
$ws.Range("H65").Value = 4088321
$ws.Range("I65").Value = 7391.225
$ws.Range("K65").Value = 36956.125
$ws.Range("M65").Value = -33836.125
$ws.Range("H74").Value = 13999.6
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13999.6
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 13999.6
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -15871.6
$ws.Range("H77").Value = 13999.6
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13999.6
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 41998.8
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -51358.8
$ws.Range("H81").Value = 105000
$ws.Range("J81").Value = 105000
$ws.Range("L81").Value = 210000
$ws.Range("N81").Value = -212122
$ws.Range("H84").Value = 105000
$ws.Range("J84").Value = 105000
$ws.Range("L84").Value = 1050000
$ws.Range("N84").Value = -1060608
$ws.Range("H94").Value = 69999
$ws.Range("J94").Value = 69999
$ws.Range("L94").Value = 69999
$ws.Range("N94").Value = -71801
$ws.Range("H117").Value = 88750
$ws.Range("J117").Value = 88750
$ws.Range("L117").Value = 88750
$ws.Range("N117").Value = -97928
$ws.Range("H136").Value = 10242
$ws.Range("I136").Value = 714.6
$ws.Range("K136").Value = 2143.8
$ws.Range("M136").Value = 406.1999999999998

